# Insert a new weekly price record as row 164 (Hortaliza / Berenjena sheet).
# Existing rows 164-179 shift down to 165-180; dimension grows from R179 to R180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(164).Insert()

$ws.Range("A164").Value = 11
$ws.Range("B164").Value = "Vega Monumental Concepción"
$ws.Range("C164").Value = "Bíobío"
$ws.Range("D164").Value = 45142
$ws.Range("D164").NumberFormat = $ws.Range("D165").NumberFormat
$ws.Range("E164").Value = 8
$ws.Range("F164").Value = 100112001
$ws.Range("G164").Value = "Berenjena"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 150
$ws.Range("K164").Value = 9000
$ws.Range("L164").Value = 9000
$ws.Range("M164").Value = 9000
$ws.Range("N164").Value = "$/caja 50 unidades"
$ws.Range("O164").Value = "Región de Arica y Parinacota"
$ws.Range("P164").Value = 180
$ws.Range("Q164").Value = 50
$ws.Range("R164").Value = "Hortaliza"
